$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.782.65"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.375.49"
$ws.Range("E3").Value = "  -3.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.92"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.18"
$ws.Range("E6").Value = "  -2.75%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  -8.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.373.21"
$ws.Range("E9").Value = "  -3.17%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("E13").Value = "  -2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.53"
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.805.03"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.613.72"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.375.01"
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("E19").Value = "  -3.49%  "
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "316.67"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.70"
$ws.Range("E22").Value = "  -3.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.83"
$ws.Range("E24").Value = "  +3.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.89"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.494.37"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("E28").Value = "  -3.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.78"
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "520.20"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.42"
$ws.Range("E31").Value = "  -4.22%  "
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("E33").Value = "  -3.27%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.66"
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.46"
$ws.Range("E38").Value = "  -6.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.374"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.99"
$ws.Range("E40").Value = "  -2.32%  "
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.35"
$ws.Range("E43").Value = "  -5.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.33"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "139.74"
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.55"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.32"
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0518"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("E50").Value = "  -1.38%  "
$ws.Range("E51").Value = "  -2.92%  "
